$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-13 down to 12-14.
$ws.Rows.Item(11).Insert()

# Copy the date style (s="2") used by column D onto the new D11 cell by
# mirroring the format of D12 (the row that used to be row 11).
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat

# Fill in the data for the newly inserted row 11.
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44855
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100114007
$ws.Cells.Item(11, 7).Value = "Jengibre"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 500
$ws.Cells.Item(11, 11).Value = 10000
$ws.Cells.Item(11, 12).Value = 10000
$ws.Cells.Item(11, 13).Value = 10000
$ws.Cells.Item(11, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 769
$ws.Cells.Item(11, 17).Value = 13
$ws.Cells.Item(11, 18).Value = "Hortaliza"
